# Update the GPS "Places" log (Sheet1) with newly logged locations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 previously held "unna" / its coordinates - replace with brugge.
$ws.Range("A2").Value = "brugge"
$ws.Range("B2").Value = "51.2089,3.2242"

# New rows logged after the trip continued through Belgium/Germany.
$ws.Range("A3").Value = "gent"
$ws.Range("B3").Value = "51.05,3.7167"

$ws.Range("A4").Value = "dusseldorf"
$ws.Range("B4").Value = "51.2217,6.7762"

$ws.Range("A5").Value = "genk"
$ws.Range("B5").Value = "50.965,5.5008"

$ws.Range("A6").Value = "brussels"
$ws.Range("B6").Value = "50.8504,4.3488"

$ws.Range("A7").Value = "antwerp"
$ws.Range("B7").Value = "51.2199,4.4035"

# Widen column B so the longer coordinate strings fit, and leave the
# selection on the last entry just like the author did after typing it.
$ws.Columns.Item(2).ColumnWidth = 14
$ws.Range("A7").Select()
